$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = 0.4635416666666667
$ws.Cells.Item(3, 2).Value = 0.3281622911694511
$ws.Cells.Item(4, 2).Value = 0.4541284403669725
$ws.Cells.Item(5, 2).Value = 0.2600979192166463

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = 0.708029197080292
$ws.Cells.Item(3, 2).Value = 0.5249266862170088
$ws.Cells.Item(4, 2).Value = 0.5043859649122807
$ws.Cells.Item(5, 2).Value = 0.489010989010989
$ws.Cells.Item(6, 2).Value = 0.6666666666666666
$ws.Cells.Item(7, 2).Value = 0.4796610169491525
$ws.Cells.Item(8, 2).Value = 0.3966666666666667

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 0.5598705501618123
$ws.Cells.Item(3, 2).Value = 0.3898840885142255
$ws.Cells.Item(4, 2).Value = 0.5032051282051282
$ws.Cells.Item(5, 2).Value = 0.6779661016949152
$ws.Cells.Item(6, 2).Value = 0.5841995841995842
$ws.Cells.Item(7, 2).Value = 0.533724340175953
$ws.Cells.Item(8, 2).Value = 0.7092198581560284
$ws.Cells.Item(9, 2).Value = 0.5654450261780105
$ws.Cells.Item(10, 2).Value = 0.5964912280701754
$ws.Cells.Item(11, 2).Value = 0.5535714285714286

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 0.4901960784313725
$ws.Cells.Item(3, 2).Value = 0.4580031695721077
$ws.Cells.Item(4, 2).Value = 0.6585365853658537
$ws.Cells.Item(5, 2).Value = 0.5168195718654435
$ws.Cells.Item(6, 2).Value = 0.3777173913043478
$ws.Cells.Item(7, 2).Value = 0.393526405451448

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 0.5539358600583091

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = 0.5278688524590164
$ws.Cells.Item(3, 2).Value = 0.6537102473498233
$ws.Cells.Item(4, 2).Value = 0.5481927710843374

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 3).Value = 0.3764825793549341
$ws.Cells.Item(2, 4).Value = 0.09914367315783958
$ws.Cells.Item(2, 5).Value = 0.2600979192166463
$ws.Cells.Item(2, 6).Value = 0.3111461981812499
$ws.Cells.Item(2, 7).Value = 0.3911453657682118
$ws.Cells.Item(2, 8).Value = 0.456481746941896
$ws.Cells.Item(2, 9).Value = 0.4635416666666667
$ws.Cells.Item(3, 3).Value = 0.5384781696432938
$ws.Cells.Item(3, 4).Value = 0.1099566009749393
$ws.Cells.Item(3, 5).Value = 0.3966666666666667
$ws.Cells.Item(3, 6).Value = 0.4843360029800707
$ws.Cells.Item(3, 7).Value = 0.5043859649122807
$ws.Cells.Item(3, 8).Value = 0.5957966764418376
$ws.Cells.Item(3, 9).Value = 0.708029197080292
$ws.Cells.Item(4, 3).Value = 0.5673577333927262
$ws.Cells.Item(4, 4).Value = 0.08849665223541958
$ws.Cells.Item(4, 5).Value = 0.3898840885142255
$ws.Cells.Item(4, 6).Value = 0.5386861122748219
$ws.Cells.Item(4, 7).Value = 0.5626577881699114
$ws.Cells.Item(4, 8).Value = 0.5934183171025276
$ws.Cells.Item(4, 9).Value = 0.7092198581560284
$ws.Cells.Item(5, 3).Value = 0.4824665336650955
$ws.Cells.Item(5, 4).Value = 0.101707153746936
$ws.Cells.Item(5, 5).Value = 0.3777173913043478
$ws.Cells.Item(5, 6).Value = 0.4096455964816129
$ws.Cells.Item(5, 7).Value = 0.4740996240017401
$ws.Cells.Item(5, 8).Value = 0.5101636985069258
$ws.Cells.Item(5, 9).Value = 0.6585365853658537
$ws.Cells.Item(6, 3).Value = 0.5539358600583091
$ws.Cells.Item(6, 5).Value = 0.5539358600583091
$ws.Cells.Item(6, 6).Value = 0.5539358600583091
$ws.Cells.Item(6, 7).Value = 0.5539358600583091
$ws.Cells.Item(6, 8).Value = 0.5539358600583091
$ws.Cells.Item(6, 9).Value = 0.5539358600583091
$ws.Cells.Item(7, 3).Value = 0.5765906236310591
$ws.Cells.Item(7, 4).Value = 0.06755621872903346
$ws.Cells.Item(7, 5).Value = 0.5278688524590164
$ws.Cells.Item(7, 6).Value = 0.5380308117716769
$ws.Cells.Item(7, 7).Value = 0.5481927710843374
$ws.Cells.Item(7, 8).Value = 0.6009515092170803
$ws.Cells.Item(7, 9).Value = 0.6537102473498233
